$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The team name in cell A17 (merged A17:A18) changes from "team 5" to "Mallika".
$ws.Range("A17").Value = "Mallika"

# Reflect the post-edit selection/scroll state (user clicked A19 after the edit).
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 17
